# Commit: "bento regression 25 scripts"
# The Sample-ID Cypher query stored in cell B3 of the "startup" sheet was
# tweaked: an intermediate WITH clause was added before `distinct lp,` and
# the final RETURN was changed to RETURN DISTINCT. Everything else (row 4's
# file query, headers, styles, etc.) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$b3Text = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE tp.chemotherapy_regimen IN ["FEC (3 week cycles)"]
  WITH  distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN DISTINCT
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
'@

$ws.Range("B3").Value = $b3Text

# The saved file's active selection also moved from C3 to B3.
$ws.Range("B3").Select()
